# Update crypto price/volume figures per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.723.13"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.599.50"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.77"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.58"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("E11").Value = "  +0.75%  "
$ws.Range("D12").Value = "1.824.57"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "1.612.24"
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.03"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("E17").Value = "  -3.28%  "
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "208.49"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -0.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.15"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E21").Value = "  +0.62%  "
$ws.Range("E22").Value = "  -3.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "143.56"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("E29").Value = "  -2.00%  "
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("E31").Value = "  +0.34%  "
$ws.Range("E32").Value = "  +0.54%  "
$ws.Range("E33").Value = "  +20.12%  "
$ws.Range("D34").Value = "1.279.50"
$ws.Range("E34").Value = "  -0.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.48"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +1.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.592"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -4.32%  "
$ws.Range("E38").Value = "  -1.69%  "
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("E40").Value = "  +0.51%  "
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.776"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "62.67"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("D44").Value = "1.735.69"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.45"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("E47").Value = "  +1.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0513"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +1.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.51"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +2.47%  "
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.399"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +1.52%  "
